$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Shift the existing table down by two rows to make room for a new
#    title row + blank spacer row above it. The style that used to
#    live on the "Ticket Hyperlink" value cell (old C4, blue/underline)
#    rides along with the shift and now sits on C6 - grab a copy of it
#    before anything else overwrites it, so it can be replanted on the
#    new hyperlink cell (C7) later.
# ------------------------------------------------------------------
$ws.Rows("1:2").Insert()
$ws.Range("C6").Copy()
$ws.Range("Z1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 2. New title row (row 1) and blank spacer row (row 2).
# ------------------------------------------------------------------
$ws.Range("A1").Value = "Fixing of FMS2134 job failure"
$ws.Range("B1").Value = ""
$ws.Range("C1").Value = ""

$ws.Range("A2").Value = ""
$ws.Range("B2").Value = ""
$ws.Range("C2").Value = ""

# ------------------------------------------------------------------
# 3. Column header row (was row 1, now row 3) - unchanged content.
# ------------------------------------------------------------------
$ws.Range("A3").Value = "Section"
$ws.Range("B3").Value = "Field"
$ws.Range("C3").Value = "Value"

# ------------------------------------------------------------------
# 4. HEADER section (rows 4-8).
# ------------------------------------------------------------------
$ws.Range("A4").Value = "HEADER"
$ws.Range("B4").Value = "Practice/Account"
$ws.Range("C4").Value = "Digital Transformation / AI Lab"

$ws.Range("A5").Value = "HEADER"
$ws.Range("B5").Value = "Project Name"
$ws.Range("C5").Value = ""

$ws.Range("A6").Value = "HEADER"
$ws.Range("B6").Value = "Ticket Hyperlink"
$ws.Range("C6").ClearFormats()
$ws.Range("C6").Value = "BP-00479"

$ws.Range("A7").Value = "HEADER"
$ws.Range("B7").Value = "Start Date"
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "15-Jul-25"

$ws.Range("A8").Value = "HEADER"
$ws.Range("B8").Value = "Deadline"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "30-Sep-25"

# ------------------------------------------------------------------
# 5. BUSINESS_CASE section (rows 9-13).
# ------------------------------------------------------------------
$ws.Range("A9").Value = "BUSINESS_CASE"
$ws.Range("B9").Value = "Why now"
$ws.Range("C9").Value = "The current manual process is slow and error-prone."

$ws.Range("A10").Value = "BUSINESS_CASE"
$ws.Range("B10").Value = "Consequences of delay"
$ws.Range("C10").Value = "Increased turnaround time and potential data entry errors."

$ws.Range("A11").Value = "BUSINESS_CASE"
$ws.Range("B11").Value = "Technical justification"
$ws.Range("C11").Value = "Leveraging LLMs for semantic validation of unstructured text."

$ws.Range("A12").Value = "BUSINESS_CASE"
$ws.Range("B12").Value = "Softtek Big Y"
$ws.Range("C12").Value = "Reduction in operational overhead."

$ws.Range("A13").Value = "BUSINESS_CASE"
$ws.Range("B13").Value = "Organizational KPI"
$ws.Range("C13").Value = "Efficiency improvement by 40%"

# ------------------------------------------------------------------
# 6. PROBLEM_STATEMENT section (rows 14-17).
# ------------------------------------------------------------------
$ws.Range("A14").Value = "PROBLEM_STATEMENT"
$ws.Range("B14").Value = "Problem Definition"
$ws.Range("C14").Value = "Project intake documents are often incomplete or inconsistent."

$ws.Range("A15").Value = "PROBLEM_STATEMENT"
$ws.Range("B15").Value = "Current Pain Points"
$ws.Range("C15").Value = "Reviewers spend hours manually checking for mandatory information."

$ws.Range("A16").Value = "PROBLEM_STATEMENT"
$ws.Range("B16").Value = "Business/System Impact"
$ws.Range("C16").Value = "Delays in project kickoff and resource allocation."

$ws.Range("A17").Value = "PROBLEM_STATEMENT"
$ws.Range("B17").Value = "Who is affected"
$ws.Range("C17").Value = "Project Management Office (PMO) and Delivery Teams."

# ------------------------------------------------------------------
# 7. PROJECT_SCOPE section (rows 18-19).
# ------------------------------------------------------------------
$ws.Range("A18").Value = "PROJECT_SCOPE"
$ws.Range("B18").Value = "In Scope"
$ws.Range("C18").Value = "Azure OpenAI integration, PDF/Excel support, LangGraph orchestration."

$ws.Range("A19").Value = "PROJECT_SCOPE"
$ws.Range("B19").Value = "Out of Scope"
$ws.Range("C19").Value = "Legacy system migration, SAP integration."

# ------------------------------------------------------------------
# 8. EXPECTED_BENEFITS section (rows 20-24).
# ------------------------------------------------------------------
$ws.Range("A20").Value = "EXPECTED_BENEFITS"
$ws.Range("B20").Value = "Qualitative Benefits"
$ws.Range("C20").Value = "Improved data quality and faster approval cycles."

$ws.Range("A21").Value = "EXPECTED_BENEFITS"
$ws.Range("B21").Value = "Softtek Hard Dollars"
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "$50,000"

$ws.Range("A22").Value = "EXPECTED_BENEFITS"
$ws.Range("B22").Value = "Softtek Soft Dollars"
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "$10,000"

$ws.Range("A23").Value = "EXPECTED_BENEFITS"
$ws.Range("B23").Value = "Customer Hard Dollars"
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "$100,000"

$ws.Range("A24").Value = "EXPECTED_BENEFITS"
$ws.Range("B24").Value = "Customer Soft Dollars"
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = "$20,000"

# ------------------------------------------------------------------
# 9. Move the ticket hyperlink from its old anchor (now C6, after the
#    shift) to the new "Start Date" cell C7, and restore the
#    blue/underline formatting (stashed in Z1 in step 1) onto C7.
# ------------------------------------------------------------------
$ws.Range("C6").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C7"), "https://jira.example.com/browse/BP-00479")

$ws.Range("Z1").Copy()
$ws.Range("C7").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("Z1").Clear()

# ------------------------------------------------------------------
# 10. Selection bookkeeping.
# ------------------------------------------------------------------
$ws.Range("A1").Select()
